$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r, date, B, C, D, E, F, G
$rows = @(
    @(149, 45064, "Krzysiu", "Kuba", "Wojtek", "Szymon", 7, 8),
    @(150, 45064, "Szymon", "Kuba", "Krzysiu", "Wojtek", 8, 2),
    @(151, 45064, "Szymon", "Krzysiu", "Kuba", "Wojtek", 5, 8),
    @(152, 45064, "Wojtek", "Szymon", "Kuba", "Krzysiu", 8, 6),
    @(153, 45064, "Szymon", "Kuba", "Eryk", "Krzysiu", 8, 4),
    @(154, 45064, "Kuba", "Krzysiu", "Wojtek", "Eryk", 8, 4),
    @(155, 45064, "Eryk", "Kuba", "Wojtek", "Szymon", 7, 8),
    @(156, 45064, "Krzysiu", "Wojtek", "Szymon", "Eryk", 4, 8),
    @(157, 45064, "Krzysiu", "Szymon", "Wojtek", "Kuba", 4, 8),
    @(158, 45064, "Eryk", "Wojtek", "Krzysiu", "Kuba", 8, 5),
    @(159, 45064, "Szymon", "Eryk", "Kuba", "Wojtek", 8, 4),
    @(160, 45064, "Szymon", "Krzysiu", "Wojtek", "Eryk", 8, 3),
    @(161, 45064, "Wojtek", "Kuba", "Eryk", "Krzysiu", 8, 5),
    @(162, 45069, "Szymon", "Kuba", "Krzysiu", "Wojtek", 8, 2),
    @(163, 45069, "Wojtek", "Szymon", "Kuba", "Krzysiu", 8, 7),
    @(164, 45069, "Krzysiu", "Szymon", "Wojtek", "Kuba", 8, 6),
    @(165, 45069, "Kuba", "Szymon", "Wojtek", "Krzysiu", 8, 1),
    @(166, 45069, "Krzysiu", "Kuba", "Szymon", "Wojtek", 3, 8),
    @(167, 45069, "Kuba", "Wojtek", "Szymon", "Krzysiu", 6, 8),
    @(168, 45069, "Szymon", "Wojtek", "Kuba", "Krzysiu", 8, 3),
    @(169, 45069, "Krzysiu", "Wojtek", "Szymon", "Kuba", 3, 8),
    @(170, 45069, "Wojtek", "Kuba", "Krzysiu", "Szymon", 5, 8),
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Formula = "=IF(F" + $r + ">G" + $r + ",1,2)"
}

$ws.Range("M159").Select() | Out-Null
